$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Duval",     215,  109,  "Garza"),
    @("Zapata",    214,  116,  "Garza"),
    @("Jim Hogg",  49,   37,   "Garza"),
    @("Starr",     340,  150,  "Garza"),
    @("McMullen",  182,  86,   "Furman"),
    @("Atascosa",  4680, 1951, "Garza"),
    @("Bexar",     7664, 3183, "Furman"),
    @("Guadalupe", 9493, 5554, "Furman")
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$ws.Range("D11").Select() | Out-Null
